$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.257.04'
$ws.Range("E2").Value = '  +0.76%  '

# Row 3
$ws.Range("D3").Value = '1.795.91'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.60'
$ws.Range("E5").Value = '  -2.25%  '

# Row 6
$ws.Range("E6").Value = '  +0.19%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4529'
$ws.Range("E7").Value = '  +16.24%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3749'
$ws.Range("E8").Value = '  +10.20%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.79'
$ws.Range("E9").Value = '  -0.99%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.148'
$ws.Range("E10").Value = '  +2.04%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07539'
$ws.Range("E11").Value = '  +4.35%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.59'
$ws.Range("E12").Value = '  +1.20%  '

# Row 13
$ws.Range("E13").Value = '  +0.22%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.303'
$ws.Range("E14").Value = '  +2.47%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.578'
$ws.Range("E15").Value = '  +7.26%  '

# Row 16
$ws.Range("D16").Value = '1.790.67'
$ws.Range("E16").Value = '  +2.02%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001092'
$ws.Range("E17").Value = '  +3.26%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06756'
$ws.Range("E18").Value = '  +2.17%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.05'
$ws.Range("E19").Value = '  +0.59%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.29%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.53'
$ws.Range("E21").Value = '  +3.53%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.348'
$ws.Range("E22").Value = '  +2.21%  '

# Row 23
$ws.Range("D23").Value = '28.224.70'
$ws.Range("E23").Value = '  +0.73%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.82'
$ws.Range("E24").Value = '  +1.53%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.419'
$ws.Range("E25").Value = '  +1.57%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.58'
$ws.Range("E26").Value = '  +3.09%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.37'
$ws.Range("E27").Value = '  -2.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.354'
$ws.Range("E28").Value = '  +1.95%  '

# Row 29
$ws.Range("D29").Value = '1.997.42'
$ws.Range("E29").Value = '  +2.08%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.10'
$ws.Range("E30").Value = '  +3.13%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.235'
$ws.Range("E31").Value = '  -3.65%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.036'
$ws.Range("E32").Value = '  -1.01%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09430'
$ws.Range("E33").Value = '  +8.72%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.816'
$ws.Range("E34").Value = '  -0.15%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2348'
$ws.Range("E35").Value = '  +11.32%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.15'

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06328'
$ws.Range("E37").Value = '  +2.83%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02333'
$ws.Range("E38").Value = '  +2.24%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.204'
$ws.Range("E39").Value = '  +1.35%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6571'
$ws.Range("E40").Value = '  +1.37%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.349'
$ws.Range("E41").Value = '  +6.18%  '

# Row 42
$ws.Range("E42").Value = '  -1.31%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.207'
$ws.Range("E43").Value = '  +0.32%  '

# Row 44
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.18%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.10'
$ws.Range("E45").Value = '  +3.04%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6106'
$ws.Range("E46").Value = '  +2.09%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.792'
$ws.Range("E47").Value = '  -0.99%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.05'
$ws.Range("E48").Value = '  +2.64%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.028'
$ws.Range("E49").Value = '  +2.37%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07122'
$ws.Range("E50").Value = '  +1.79%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.164'
$ws.Range("E51").Value = '  +1.04%  '
